$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table ("テーブル1" / table1.xml) currently spans C3:E6 with 3 columns
# (sut, prefix, expected). Extend it by one column (dummy) to C3:F6.
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("C3:F6"))

# Header for the newly-added 4th column.
$ws.Range("F3").Value = "dummy"

# Body rows for the dummy column - all literal empty-quote strings.
$ws.Range("F4").Value = '""'
$ws.Range("F5").Value = '""'
$ws.Range("F6").Value = '""'

# Move the active selection the same way Excel would after typing down
# the new column (was E7, now F7).
[void]$ws.Range("F7").Select()
